$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = [double]"3"
$ws.Range("F2").Value2 = [double]"1"
$ws.Range("G2").Value2 = [double]"0.7169949999999999"
$ws.Range("H2").Value2 = [double]"2.150985"
$ws.Range("I2").Value2 = [double]"0.003651663653539308"
$ws.Range("J2").Value2 = [double]"0.003651663653539308"
$ws.Range("M2").Value2 = [double]"194.2587713333333"
$ws.Range("N2").Value2 = [double]"582.776314"
$ws.Range("O2").Value2 = [double]"0.9908911906753961"
$ws.Range("P2").Value2 = [double]"0.9908911906753961"
$ws.Range("Q2").Value2 = [double]"139.2825677521433"
$ws.Range("R2").Value2 = [double]"1253.54310976929"
$ws.Range("S2").Value2 = [double]"0.003618401345601632"
$ws.Range("T2").Value2 = [double]"0.003618401345601632"
$ws.Range("E3").Value2 = [double]"3"
$ws.Range("F3").Value2 = [double]"1"
$ws.Range("G3").Value2 = [double]"0.7169949999999999"
$ws.Range("H3").Value2 = [double]"2.150985"
$ws.Range("I3").Value2 = [double]"0.003651663653539308"
$ws.Range("J3").Value2 = [double]"0.003651663653539308"
$ws.Range("O3").Value2 = [double]"0.004009241031003317"
$ws.Range("P3").Value2 = [double]"0.004009241031003317"
$ws.Range("Q3").Value2 = [double]"0.5635506610516665"
$ws.Range("R3").Value2 = [double]"5.071955949464999"
$ws.Range("S3").Value2 = [double]"1.464039975119328E-05"
$ws.Range("T3").Value2 = [double]"1.464039975119328E-05"
$ws.Range("E4").Value2 = [double]"3"
$ws.Range("F4").Value2 = [double]"1"
$ws.Range("G4").Value2 = [double]"0.7169949999999999"
$ws.Range("H4").Value2 = [double]"2.150985"
$ws.Range("I4").Value2 = [double]"0.003651663653539308"
$ws.Range("J4").Value2 = [double]"0.003651663653539308"
$ws.Range("M4").Value2 = [double]"0.9314810000000001"
$ws.Range("N4").Value2 = [double]"2.794443"
$ws.Range("O4").Value2 = [double]"0.004751375244712719"
$ws.Range("P4").Value2 = [double]"0.004751375244712718"
$ws.Range("Q4").Value2 = [double]"0.6678672195950001"
$ws.Range("R4").Value2 = [double]"6.010804976355001"
$ws.Range("S4").Value2 = [double]"1.735042428544387E-05"
$ws.Range("T4").Value2 = [double]"1.735042428544387E-05"
$ws.Range("E5").Value2 = [double]"3"
$ws.Range("F5").Value2 = [double]"1"
$ws.Range("G5").Value2 = [double]"0.7169949999999999"
$ws.Range("H5").Value2 = [double]"2.150985"
$ws.Range("I5").Value2 = [double]"0.003651663653539308"
$ws.Range("J5").Value2 = [double]"0.003651663653539308"
$ws.Range("M5").Value2 = [double]"0.06826133333333334"
$ws.Range("N5").Value2 = [double]"0.204784"
$ws.Range("O5").Value2 = [double]"0.0003481930488878283"
$ws.Range("P5").Value2 = [double]"0.0003481930488878283"
$ws.Range("Q5").Value2 = [double]"0.04894303469333333"
$ws.Range("R5").Value2 = [double]"0.44048731224"
$ws.Range("S5").Value2 = [double]"1.271483901038718E-06"
$ws.Range("T5").Value2 = [double]"1.271483901038718E-06"
$ws.Range("I6").Value2 = [double]"0.9751961860217362"
$ws.Range("J6").Value2 = [double]"0.9751961860217361"
$ws.Range("M6").Value2 = [double]"194.2587713333333"
$ws.Range("N6").Value2 = [double]"582.776314"
$ws.Range("O6").Value2 = [double]"0.9908911906753961"
$ws.Range("P6").Value2 = [double]"0.9908911906753961"
$ws.Range("Q6").Value2 = [double]"37196.14995744627"
$ws.Range("R6").Value2 = [double]"334765.3496170164"
$ws.Range("S6").Value2 = [double]"0.9663133099091833"
$ws.Range("T6").Value2 = [double]"0.9663133099091832"
$ws.Range("I7").Value2 = [double]"0.9751961860217362"
$ws.Range("J7").Value2 = [double]"0.9751961860217361"
$ws.Range("O7").Value2 = [double]"0.004009241031003317"
$ws.Range("P7").Value2 = [double]"0.004009241031003317"
$ws.Range("S7").Value2 = [double]"0.003909796562276288"
$ws.Range("T7").Value2 = [double]"0.003909796562276288"
$ws.Range("I8").Value2 = [double]"0.9751961860217362"
$ws.Range("J8").Value2 = [double]"0.9751961860217361"
$ws.Range("M8").Value2 = [double]"0.9314810000000001"
$ws.Range("N8").Value2 = [double]"2.794443"
$ws.Range("O8").Value2 = [double]"0.004751375244712719"
$ws.Range("P8").Value2 = [double]"0.004751375244712718"
$ws.Range("Q8").Value2 = [double]"178.3574904788187"
$ws.Range("R8").Value2 = [double]"1605.217414309368"
$ws.Range("S8").Value2 = [double]"0.004633523017001937"
$ws.Range("T8").Value2 = [double]"0.004633523017001936"
$ws.Range("I9").Value2 = [double]"0.9751961860217362"
$ws.Range("J9").Value2 = [double]"0.9751961860217361"
$ws.Range("M9").Value2 = [double]"0.06826133333333334"
$ws.Range("N9").Value2 = [double]"0.204784"
$ws.Range("O9").Value2 = [double]"0.0003481930488878283"
$ws.Range("P9").Value2 = [double]"0.0003481930488878283"
$ws.Range("Q9").Value2 = [double]"13.07049753035378"
$ws.Range("R9").Value2 = [double]"117.634477773184"
$ws.Range("S9").Value2 = [double]"0.0003395565332746901"
$ws.Range("T9").Value2 = [double]"0.00033955653327469"
$ws.Range("G10").Value2 = [double]"4.138615666666666"
$ws.Range("H10").Value2 = [double]"12.415847"
$ws.Range("I10").Value2 = [double]"0.02107801645190694"
$ws.Range("J10").Value2 = [double]"0.02107801645190694"
$ws.Range("M10").Value2 = [double]"194.2587713333333"
$ws.Range("N10").Value2 = [double]"582.776314"
$ws.Range("O10").Value2 = [double]"0.9908911906753961"
$ws.Range("P10").Value2 = [double]"0.9908911906753961"
$ws.Range("Q10").Value2 = [double]"803.9623944275509"
$ws.Range("R10").Value2 = [double]"7235.661549847957"
$ws.Range("S10").Value2 = [double]"0.02088602081910566"
$ws.Range("T10").Value2 = [double]"0.02088602081910566"
$ws.Range("G11").Value2 = [double]"4.138615666666666"
$ws.Range("H11").Value2 = [double]"12.415847"
$ws.Range("I11").Value2 = [double]"0.02107801645190694"
$ws.Range("J11").Value2 = [double]"0.02107801645190694"
$ws.Range("O11").Value2 = [double]"0.004009241031003317"
$ws.Range("P11").Value2 = [double]"0.004009241031003317"
$ws.Range("Q11").Value2 = [double]"3.252909148304777"
$ws.Range("R11").Value2 = [double]"29.27618233474299"
$ws.Range("S11").Value2 = [double]"8.450684841114828E-05"
$ws.Range("T11").Value2 = [double]"8.450684841114827E-05"
$ws.Range("G12").Value2 = [double]"4.138615666666666"
$ws.Range("H12").Value2 = [double]"12.415847"
$ws.Range("I12").Value2 = [double]"0.02107801645190694"
$ws.Range("J12").Value2 = [double]"0.02107801645190694"
$ws.Range("M12").Value2 = [double]"0.9314810000000001"
$ws.Range("N12").Value2 = [double]"2.794443"
$ws.Range("O12").Value2 = [double]"0.004751375244712719"
$ws.Range("P12").Value2 = [double]"0.004751375244712718"
$ws.Range("Q12").Value2 = [double]"3.855041859802334"
$ws.Range("R12").Value2 = [double]"34.695376738221"
$ws.Range("S12").Value2 = [double]"0.0001001495655772381"
$ws.Range("T12").Value2 = [double]"0.000100149565577238"
$ws.Range("G13").Value2 = [double]"4.138615666666666"
$ws.Range("H13").Value2 = [double]"12.415847"
$ws.Range("I13").Value2 = [double]"0.02107801645190694"
$ws.Range("J13").Value2 = [double]"0.02107801645190694"
$ws.Range("M13").Value2 = [double]"0.06826133333333334"
$ws.Range("N13").Value2 = [double]"0.204784"
$ws.Range("O13").Value2 = [double]"0.0003481930488878283"
$ws.Range("P13").Value2 = [double]"0.0003481930488878283"
$ws.Range("Q13").Value2 = [double]"0.2825074235608889"
$ws.Range("R13").Value2 = [double]"2.542566812048"
$ws.Range("S13").Value2 = [double]"7.339218812897282E-06"
$ws.Range("T13").Value2 = [double]"7.339218812897282E-06"
$ws.Range("E14").Value2 = [double]"1"
$ws.Range("F14").Value2 = [double]"0.3333333333333333"
$ws.Range("G14").Value2 = [double]"0.014556"
$ws.Range("H14").Value2 = [double]"0.043668"
$ws.Range("I14").Value2 = [double]"7.413387281768795E-05"
$ws.Range("J14").Value2 = [double]"7.413387281768795E-05"
$ws.Range("M14").Value2 = [double]"194.2587713333333"
$ws.Range("N14").Value2 = [double]"582.776314"
$ws.Range("O14").Value2 = [double]"0.9908911906753961"
$ws.Range("P14").Value2 = [double]"0.9908911906753961"
$ws.Range("Q14").Value2 = [double]"2.827630675528"
$ws.Range("R14").Value2 = [double]"25.448676079752"
$ws.Range("S14").Value2 = [double]"7.345860150569719E-05"
$ws.Range("T14").Value2 = [double]"7.345860150569719E-05"
$ws.Range("E15").Value2 = [double]"1"
$ws.Range("F15").Value2 = [double]"0.3333333333333333"
$ws.Range("G15").Value2 = [double]"0.014556"
$ws.Range("H15").Value2 = [double]"0.043668"
$ws.Range("I15").Value2 = [double]"7.413387281768795E-05"
$ws.Range("J15").Value2 = [double]"7.413387281768795E-05"
$ws.Range("O15").Value2 = [double]"0.004009241031003317"
$ws.Range("P15").Value2 = [double]"0.004009241031003317"
$ws.Range("Q15").Value2 = [double]"0.011440865588"
$ws.Range("R15").Value2 = [double]"0.102967790292"
$ws.Range("S15").Value2 = [double]"2.97220564687856E-07"
$ws.Range("T15").Value2 = [double]"2.97220564687856E-07"
$ws.Range("E16").Value2 = [double]"1"
$ws.Range("F16").Value2 = [double]"0.3333333333333333"
$ws.Range("G16").Value2 = [double]"0.014556"
$ws.Range("H16").Value2 = [double]"0.043668"
$ws.Range("I16").Value2 = [double]"7.413387281768795E-05"
$ws.Range("J16").Value2 = [double]"7.413387281768795E-05"
$ws.Range("M16").Value2 = [double]"0.9314810000000001"
$ws.Range("N16").Value2 = [double]"2.794443"
$ws.Range("O16").Value2 = [double]"0.004751375244712719"
$ws.Range("P16").Value2 = [double]"0.004751375244712718"
$ws.Range("Q16").Value2 = [double]"0.013558637436"
$ws.Range("R16").Value2 = [double]"0.122027736924"
$ws.Range("S16").Value2 = [double]"3.522378481006437E-07"
$ws.Range("T16").Value2 = [double]"3.522378481006436E-07"
$ws.Range("E17").Value2 = [double]"1"
$ws.Range("F17").Value2 = [double]"0.3333333333333333"
$ws.Range("G17").Value2 = [double]"0.014556"
$ws.Range("H17").Value2 = [double]"0.043668"
$ws.Range("I17").Value2 = [double]"7.413387281768795E-05"
$ws.Range("J17").Value2 = [double]"7.413387281768795E-05"
$ws.Range("M17").Value2 = [double]"0.06826133333333334"
$ws.Range("N17").Value2 = [double]"0.204784"
$ws.Range("O17").Value2 = [double]"0.0003481930488878283"
$ws.Range("P17").Value2 = [double]"0.0003481930488878283"
$ws.Range("Q17").Value2 = [double]"0.000993611968"
$ws.Range("R17").Value2 = [double]"0.008942507712000001"
$ws.Range("S17").Value2 = [double]"2.581289920225326E-08"
$ws.Range("T17").Value2 = [double]"2.581289920225326E-08"
